# Insert a new "Predicted Salary" column between "Current Year Salary" (B)
# and "Residual" (C), shifting the existing Residual column to D, then
# populate it with each player's predicted salary (Current Year Salary -
# Residual), computed from the existing data already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift column C (Residual) one place to the right, freeing up C for the
# new "Predicted Salary" column. Insert() on EntireColumn carries over the
# formatting (bold/centered header style) from the column it pushes out of
# the way, so the new header cell already matches the other header cells.
$ws.Range("C1").EntireColumn.Insert()

# Header
$ws.Range("C1").Value = "Predicted Salary"

# Find the last used row so we cover every data row regardless of count.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $salary = $ws.Cells.Item($r, 2).Value2
    $residual = $ws.Cells.Item($r, 4).Value2
    if ($salary -ne $null -and $residual -ne $null) {
        $ws.Cells.Item($r, 3).Value = $salary - $residual
    }
}
